$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" updates ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E3").Value = 771.7913970000002
$schedule.Range("F3").Value = 29.16823117913833
$schedule.Range("E4").Value = 371.639775
$schedule.Range("F4").Value = 10.92415564373898

# --- Sheet "Detailed" updates ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B47").Value = 56.98
$detailed.Range("B48").Value = 57.03878

$detailed.Range("B49").Value = 57.06
$detailed.Range("C49").Value = "historical"

$detailed.Range("B53").Value = 56.98
$detailed.Range("B54").Value = 56.97994

$detailed.Range("B61").Value = 58.43713
$detailed.Range("B62").Value = 58.00918

$detailed.Range("B64").Value = 26.82191
$detailed.Range("B65").Value = 34.69723
$detailed.Range("B66").Value = 25.178
$detailed.Range("B67").Value = 32.18984

$detailed.Range("B70").Value = 0.02957

$detailed.Range("B72").Value = 5.26268
$detailed.Range("B73").Value = 22.07
$detailed.Range("B74").Value = 23.31255
$detailed.Range("B75").Value = 33.78973
$detailed.Range("B76").Value = 34.45564

$detailed.Range("B80").Value = 33.2633
$detailed.Range("B81").Value = 0.00855
$detailed.Range("B82").Value = -1.16441
$detailed.Range("B83").Value = -4.62815
$detailed.Range("B84").Value = -4.57905
$detailed.Range("B85").Value = 4.39442
$detailed.Range("B86").Value = 35.14435
$detailed.Range("B87").Value = 45.7518
$detailed.Range("B88").Value = 57.03529
$detailed.Range("B89").Value = 58.40626
$detailed.Range("B90").Value = 58.48021
$detailed.Range("B91").Value = 57.06

$detailed.Range("B93").Value = 57.06
